$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RuntimesChart")

# Add the new data point for day 3 of Advent of Code 2025
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 0.00095804

# Update selection to mirror the author's final state
$ws.Range("A5:B5").Select()
